$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'31.017.47"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  +1.17%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'1.958.34"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  -0.28%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('D4').Value = "'1.001"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = "'  +0.11%  "
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'245.06"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  -2.01%  "
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'0.9989"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  -0.10%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = "'0.4876"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = "'  +0.74%  "
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Value = "'0.2960"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = "'  +0.32%  "
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = "'0.06818"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'  +0.40%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = "'19.18"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'  -1.24%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = "'107.06"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'  -2.77%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range('B12').Value = "'TRON"
$ws.Range('B12').Style = 'Normal'
$ws.Range('C12').Value = "'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range('C12').Style = 'Normal'
$ws.Range('D12').Value = "'0.07822"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'  +0.66%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('B13').Value = "'WrappedEther"
$ws.Range('B13').Style = 'Normal'
$ws.Range('C13').Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range('C13').Style = 'Normal'
$ws.Range('D13').Value = "'1.944.52"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'  -1.06%  "
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'5.491"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'  +0.45%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'0.7036"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  +2.08%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'284.40"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'  -3.54%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = "'31.014.12"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'  +1.09%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range('B18').Value = "'Avalanche"
$ws.Range('B18').Style = 'Normal'
$ws.Range('C18').Value = "'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range('C18').Style = 'Normal'
$ws.Range('D18').Value = "'13.20"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'  -0.42%  "
$ws.Range('E18').Style = 'Normal'
$ws.Range('B19').Value = "'WrappedliquidstakedEther2.0"
$ws.Range('B19').Style = 'Normal'
$ws.Range('C19').Value = "'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range('C19').Style = 'Normal'
$ws.Range('D19').Value = "'2.227.67"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'  -0.10%  "
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = "'0.000007692"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'  -0.28%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = "'0.9993"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'  +0.01%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'5.509"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'  -1.95%  "
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'1.001"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'  +0.13%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'6.500"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'  -1.75%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = "'9.778"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'  -1.36%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = "'168.80"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'  -1.03%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = "'19.99"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'  -1.11%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range('E28').Value = "'  +0.74%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('E29').Value = "'  -0.93%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = "'1.402"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'  -2.32%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = "'1.585"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'  -1.80%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = "'4.611"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'  -2.44%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = "'4.435"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'  -0.30%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = "'0.04922"
$ws.Range('D34').Style = 'Normal'
$ws.Range('D35').Value = "'0.7610"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "'  -1.52%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = "'1.173"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'  -0.73%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = "'2.733"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'  +0.04%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = "'0.02011"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'  -1.84%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = "'2.701"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'  -0.73%  "
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = "'6.542"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'  +2.54%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = "'77.19"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'  +9.61%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = "'2.113"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'  -0.36%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = "'0.8900"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'  +1.42%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = "'0.4470"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'  -0.27%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = "'108.89"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'  -0.36%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = "'8.114"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'  +7.83%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = "'0.9997"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'  -0.09%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = "'1.004.92"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'  +9.31%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = "'0.1259"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'  -1.98%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = "'9.338"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'  -0.33%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = "'35.93"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'  -0.58%  "
$ws.Range('E51').Style = 'Normal'
